$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.946.40'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '1.675.30'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.28%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0888'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = '1.912.33'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '1.700.56'
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '26.971.05'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '235.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("D20").Value = '0.0₃0733'
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("E24").Value = '  -4.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("E28").Value = '  -3.69%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  +0.38%  '
$ws.Range("E31").Value = '  -1.04%  '
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("D33").Value = '1.486.67'
$ws.Range("E33").Value = '  +2.21%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.68'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.50%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.585'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.07%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0175'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.63%  '
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.78%  '
$ws.Range("E41").Value = '  +4.58%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '67.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.08%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("D45").Value = '1.819.67'
$ws.Range("E45").Value = '  +0.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.778'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("E50").Value = '  +2.21%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0508'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.03%  '
